# Applies the OOXML changes described in the commit:
#   "Added output to Excel file function"
#
# Functional changes on the "Master Scores" sheet:
#   - B2 changes from -1000 to 0
#   - C2 changes from -1000 to 0
#   - F20 changes from 6 to 100000
#   - The active selection moves from C3 to F20 (the last cell touched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zero-out the previous "penalty" values in row 2.
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# Write the new "output" value into F20.
$ws.Range("F20").Value = 100000

# Leave the selection on the last cell written, matching the saved
# workbook view (activeCell/sqref = F20) captured when the author's
# script finished running.
$ws.Range("F20").Select()

# Best-effort: also mirror the workbook window chrome (size/position)
# recorded in the diff. This is cosmetic window state that some hosts
# don't persist, so failures here are harmless and ignored.
try {
    $win = $excel.ActiveWindow
    $win.Left = 0
    $win.Top = 0
    $win.Width = 25600
    $win.Height = 16060
} catch {}
